# Applies the TempBasket / OrderHistory data changes captured in the commit.
$wb = $excel.ActiveWorkbook

# --- TempBasket sheet -------------------------------------------------
$basket = $wb.Worksheets.Item("TempBasket")

# Row 8-10 used to hold CafeMate0007's basket (Espresso/Wrap/DonkAsk/Three
# with quantities). That basket is now recorded differently: one item
# (Fires) spread across rows 8-10, one id+value per row. Clear the old
# C:E values first since the new layout only uses columns A and B.
$basket.Range("C8:E10").ClearContents()

$basket.Range("A8").Value = "CafeMate0007"
$basket.Range("B8").Value = "Fires"

$basket.Range("A9").Value = "CafeMate0007"
$basket.Range("B9").Value = 1

$basket.Range("A10").Value = "CafeMate0007"
$basket.Range("B10").Value = 1

# Row 11 gains admin's new basket header (SoftDrink/Coke/ChickenWrap/Three)
$basket.Range("A11").Value = "admin"
$basket.Range("B11").Value = "SoftDrink"
$basket.Range("C11").Value = "Coke"
$basket.Range("D11").Value = "ChickenWrap"
$basket.Range("E11").Value = "Three"

# Row 12 gains the matching quantities
$basket.Range("A12").Value = "admin"
$basket.Range("B12").Value = 1
$basket.Range("C12").Value = 1
$basket.Range("D12").Value = 1
$basket.Range("E12").Value = 1

# Row 13 gains another copy of the quantities
$basket.Range("A13").Value = "admin"
$basket.Range("B13").Value = 1
$basket.Range("C13").Value = 1
$basket.Range("D13").Value = 1
$basket.Range("E13").Value = 1

# --- OrderHistory sheet -------------------------------------------------
$history = $wb.Worksheets.Item("OrderHistory")

$history.Range("A6").Value = "CafeMate0007"
$history.Range("B6").Value = "Basket CafeMate0007 has the following items:`nItem: Espresso, price: 1.00 pounds`nItem: Wrap, price: 1.00 pounds`nItem: DonkAsk, price: 1.00 pounds`nItem: Three, price: 1.00 pounds, quantity: 4`nTotal cost: 7.00 pounds`n"

$history.Range("A7").Value = "CafeMate0007"
$history.Range("B7").Value = "Basket CafeMate0007 has the following items:`nItem: SoftDrink, price: 1.00 pounds`nTotal cost: 1.00 pounds`n"
